$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (70) with the daily update values, matching the
# existing date formatting used in column A.
$ws.Range("A70").Value = 46019
$ws.Range("A70").NumberFormat = $ws.Range("A69").NumberFormat
$ws.Range("B70").Value = 154
$ws.Range("C70").Value = 164
$ws.Range("D70").Value = 154
